# feat: add 2022-Q1 data
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert a brand-new sheet "2022-Q1" right before the "总计" sheet,
#    holding the per-fund breakdown for the new quarter.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Re-resolve both sheets by name now that the tab order changed - a
# cached reference can silently point at "whatever sheet now sits at
# that index" instead of the sheet it originally pointed to.
$newSheet = $wb.Worksheets.Item("2022-Q1")
$prevQuarter = $wb.Worksheets.Item("2021-Q4")

# Pull over the header / index-column formatting from the previous
# quarter's sheet so the new sheet matches the established look.
$prevQuarter.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$prevQuarter.Range("A2").Copy()
$newSheet.Range("A2:A5").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows. Columns B-G hold text values in the source data (fund
# codes/names and numeric-looking figures alike are stored as text, so
# a code like "006323" keeps its leading zero) - briefly force text
# format before writing so COM doesn't silently coerce the numeric-
# looking strings to numbers, then drop the format again so the cells
# end up unstyled like the rest of the data rows.
$dataText = $newSheet.Range("B2:G5")
$dataText.NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "006323"
$newSheet.Range("C2").Value = "合煦智远嘉选混合A"
$newSheet.Range("D2").Value = "1.54"
$newSheet.Range("E2").Value = "79.45"
$newSheet.Range("F2").Value = "5.60"
$newSheet.Range("G2").Value = "0.0862"
$newSheet.Range("H2").Value = 3

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "006324"
$newSheet.Range("C3").Value = "合煦智远嘉选混合C"
$newSheet.Range("D3").Value = "0.59"
$newSheet.Range("E3").Value = "79.45"
$newSheet.Range("F3").Value = "5.60"
$newSheet.Range("G3").Value = "0.0330"
$newSheet.Range("H3").Value = 3

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "010783"
$newSheet.Range("C4").Value = "德邦沪港深龙头混合A"
$newSheet.Range("D4").Value = "0.93"
$newSheet.Range("E4").Value = "81.58"
$newSheet.Range("F4").Value = "2.66"
$newSheet.Range("G4").Value = "0.0247"
$newSheet.Range("H4").Value = 10

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "010784"
$newSheet.Range("C5").Value = "德邦沪港深龙头混合C"
$newSheet.Range("D5").Value = "0.27"
$newSheet.Range("E5").Value = "81.58"
$newSheet.Range("F5").Value = "2.66"
$newSheet.Range("G5").Value = "0.0072"
$newSheet.Range("H5").Value = 10

$dataText.ClearFormats()

# ---------------------------------------------------------------------
# 2) Add the 2022-Q1 summary row at the top of the "总计" sheet's data
#    (everything else shifts down by one row).
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# New row lands with inherited header formatting - strip it back to the
# plain (unstyled) look the other data rows use.
$totalSheet.Range("B2:D2").ClearFormats()

# The index column (A) keeps the bordered/centered "s=2" look used by
# every other row - copy that formatting down from row 3.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 0.15

# Column A is a literal 0-based row index, not a formula - renumber the
# rows that shifted down so the sequence stays 0,1,2,3,4,5.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
